$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tracked dates (A1:C1) advance by one day: 2023-05-12 -> 2023-05-13
$ws.Range("A1:C1").Value = 45059

# J1:K1:L1 were blank placeholder cells carrying the date style only;
# clearing them drops them out of the sheet entirely (no value, no style left behind).
$ws.Range("J1:L1").Clear()

# Row 2 "Work" flags: A2 flips from No (0) to Yes (100); C2 flips from No (0) to Yes (100).
$ws.Range("A2").Value = 100
$ws.Range("C2").Value = 100

# A new (currently empty) tracking row is started at row 5, with B5 pre-formatted
# like the date cells in row 1 (so it shares the same style index, no new style created).
$ws.Range("A1").Copy()
$ws.Range("B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Select the whole sheet (mirrors the "select all" gesture captured in the saved view state).
$ws.Cells.Select()
